$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header text updates (rich-text shared strings) - edit only the numeric /
# date substrings via Characters() so the surrounding run formatting and
# the rest of the sentence are left alone.
# ---------------------------------------------------------------------------

# A8: "Volume 32   Number  41" -> "...42"
$volCell = $ws.Cells.Item(8, 1)
$volCell.Characters(21, 2).Text = "42"

# C9: "Report Covering the Week  10/6/2025  Through  10/12/2025"
#   -> "...10/13/2025  Through  10/19/2025"
$weekCell = $ws.Cells.Item(9, 3)
$weekCell.Characters(27, 9).Text = "10/13/2025"
$weekCell.Characters(48, 10).Text = "10/19/2025"

# ---------------------------------------------------------------------------
# Helper reference cells used to clone formatting (style + shared string)
# when a cell needs to switch between "numeric" and "text placeholder"
# representations. These particular cells are never themselves modified
# by this script, so they stay valid as templates throughout.
# ---------------------------------------------------------------------------
$refTextZero = $ws.Cells.Item(14, 3)    # C14: s="13" t="s" -> shared string "0"
$refTextStar = $ws.Cells.Item(14, 5)    # E14: s="13" t="s" -> shared string "***.*"
$refNumCount = $ws.Cells.Item(14, 9)    # I14: s="14" plain integer style
$refNumPct   = $ws.Cells.Item(14, 11)   # K14: s="15" plain percent style

function Set-Num($cell, $value) {
    # Ensure the cell uses the plain integer/number style, then set value.
    $refNumCount.Copy($cell)
    $cell.Value = $value
}

function Set-Pct($cell, $value) {
    # Ensure the cell uses the plain percent style, then set value.
    $refNumPct.Copy($cell)
    $cell.Value = $value
}

function Set-TextZero($cell) {
    $refTextZero.Copy($cell)
}

function Set-TextStar($cell) {
    $refTextStar.Copy($cell)
}

# ---------------------------------------------------------------------------
# Row 14
# ---------------------------------------------------------------------------
$ws.Cells.Item(14, 14).Value = -87.5   # N14

# ---------------------------------------------------------------------------
# Row 15
# ---------------------------------------------------------------------------
Set-TextZero $ws.Cells.Item(15, 6)     # F15 -> text "0"
$ws.Cells.Item(15, 7).Value = 1        # G15
$ws.Cells.Item(15, 8).Value = -100     # H15

# ---------------------------------------------------------------------------
# Row 16
# ---------------------------------------------------------------------------
$ws.Cells.Item(16, 4).Value = 1                      # D16
$ws.Cells.Item(16, 5).Value = 100                     # E16
$ws.Cells.Item(16, 7).Value = 9                       # G16
$ws.Cells.Item(16, 8).Value = 66.666666666666         # H16
$ws.Cells.Item(16, 9).Value = 97                      # I16
$ws.Cells.Item(16, 10).Value = 116                    # J16
$ws.Cells.Item(16, 11).Value = -16.379310344827       # K16
$ws.Cells.Item(16, 12).Value = -22.4                  # L16
$ws.Cells.Item(16, 13).Value = -23.622047244094       # M16
$ws.Cells.Item(16, 14).Value = -79.664570230608       # N16

# ---------------------------------------------------------------------------
# Row 17
# ---------------------------------------------------------------------------
$ws.Cells.Item(17, 3).Value = 3                       # C17
$ws.Cells.Item(17, 4).Value = 4                       # D17
$ws.Cells.Item(17, 5).Value = -25                     # E17
$ws.Cells.Item(17, 7).Value = 16                      # G17
$ws.Cells.Item(17, 8).Value = -37.5                   # H17
$ws.Cells.Item(17, 9).Value = 182                     # I17
$ws.Cells.Item(17, 10).Value = 198                    # J17
$ws.Cells.Item(17, 11).Value = -8.080808080808        # K17
$ws.Cells.Item(17, 12).Value = 8.982035928143         # L17
$ws.Cells.Item(17, 13).Value = 91.578947368421        # M17
$ws.Cells.Item(17, 14).Value = -19.823788546255       # N17

# ---------------------------------------------------------------------------
# Row 18
# ---------------------------------------------------------------------------
$ws.Cells.Item(18, 3).Value = 2                       # C18
Set-TextZero $ws.Cells.Item(18, 4)                    # D18 -> text "0"
Set-TextStar $ws.Cells.Item(18, 5)                    # E18 -> text "***.*"
$ws.Cells.Item(18, 6).Value = 9                       # F18
$ws.Cells.Item(18, 7).Value = 9                       # G18
$ws.Cells.Item(18, 8).Value = 0                       # H18
$ws.Cells.Item(18, 9).Value = 110                     # I18
$ws.Cells.Item(18, 11).Value = -16.666666666666       # K18
$ws.Cells.Item(18, 12).Value = -5.982905982905        # L18
$ws.Cells.Item(18, 13).Value = -17.293233082706       # M18
$ws.Cells.Item(18, 14).Value = -88.197424892703       # N18

# ---------------------------------------------------------------------------
# Row 19
# ---------------------------------------------------------------------------
$ws.Cells.Item(19, 3).Value = 8                       # C19
$ws.Cells.Item(19, 4).Value = 13                      # D19
$ws.Cells.Item(19, 5).Value = -38.461538461538        # E19
$ws.Cells.Item(19, 6).Value = 37                      # F19
$ws.Cells.Item(19, 7).Value = 52                      # G19
$ws.Cells.Item(19, 8).Value = -28.846153846153        # H19
$ws.Cells.Item(19, 9).Value = 456                     # I19
$ws.Cells.Item(19, 10).Value = 619                    # J19
$ws.Cells.Item(19, 11).Value = -26.332794830371       # K19
$ws.Cells.Item(19, 12).Value = -16.022099447513       # L19
$ws.Cells.Item(19, 13).Value = 68.265682656826        # M19
$ws.Cells.Item(19, 14).Value = 23.913043478260        # N19

# ---------------------------------------------------------------------------
# Row 20
# ---------------------------------------------------------------------------
$ws.Cells.Item(20, 3).Value = 3                       # C20
$ws.Cells.Item(20, 4).Value = 7                       # D20
$ws.Cells.Item(20, 5).Value = -57.142857142857        # E20
$ws.Cells.Item(20, 6).Value = 18                      # F20
$ws.Cells.Item(20, 7).Value = 23                      # G20
$ws.Cells.Item(20, 8).Value = -21.739130434782        # H20
$ws.Cells.Item(20, 9).Value = 219                     # I20
$ws.Cells.Item(20, 10).Value = 270                    # J20
$ws.Cells.Item(20, 11).Value = -18.888888888888       # K20
$ws.Cells.Item(20, 12).Value = -21.223021582733       # L20
$ws.Cells.Item(20, 13).Value = 114.705882352941       # M20
$ws.Cells.Item(20, 14).Value = -85.102040816326       # N20

# ---------------------------------------------------------------------------
# Row 21
# ---------------------------------------------------------------------------
$ws.Cells.Item(21, 3).Value = 18                      # C21
$ws.Cells.Item(21, 4).Value = 25                      # D21
$ws.Cells.Item(21, 5).Value = -28                     # E21
$ws.Cells.Item(21, 6).Value = 89                      # F21
$ws.Cells.Item(21, 7).Value = 110                     # G21
$ws.Cells.Item(21, 8).Value = -19.090909090909        # H21
$ws.Cells.Item(21, 9).Value = 1086                    # I21
$ws.Cells.Item(21, 10).Value = 1347                   # J21
$ws.Cells.Item(21, 11).Value = -19.376391982182       # K21
$ws.Cells.Item(21, 12).Value = -12.630732099758       # L21
$ws.Cells.Item(21, 13).Value = 46.756756756756        # M21
$ws.Cells.Item(21, 14).Value = -68.980291345329       # N21

# ---------------------------------------------------------------------------
# Row 22 (several numeric <-> text swaps)
# ---------------------------------------------------------------------------
Set-Num $ws.Cells.Item(22, 4) 1        # D22 -> numeric 1
Set-Pct $ws.Cells.Item(22, 5) -100     # E22 -> numeric -100
Set-TextZero $ws.Cells.Item(22, 6)     # F22 -> text "0"
Set-Num $ws.Cells.Item(22, 7) 1        # G22 -> numeric 1
Set-Pct $ws.Cells.Item(22, 8) -100     # H22 -> numeric -100
$ws.Cells.Item(22, 10).Value = 5       # J22
$ws.Cells.Item(22, 11).Value = 40      # K22
$ws.Cells.Item(22, 12).Value = -22.222222222222  # L22

# ---------------------------------------------------------------------------
# Row 23
# ---------------------------------------------------------------------------
$ws.Cells.Item(23, 3).Value = 2                       # C23
Set-Num $ws.Cells.Item(23, 4) 2        # D23 -> numeric 2
Set-Pct $ws.Cells.Item(23, 5) 0        # E23 -> numeric 0
$ws.Cells.Item(23, 6).Value = 5                       # F23
$ws.Cells.Item(23, 8).Value = 150                     # H23
$ws.Cells.Item(23, 9).Value = 45                      # I23
$ws.Cells.Item(23, 10).Value = 53                     # J23
$ws.Cells.Item(23, 11).Value = -15.094339622641       # K23
$ws.Cells.Item(23, 12).Value = -31.818181818181       # L23
$ws.Cells.Item(23, 13).Value = 36.363636363636        # M23

# ---------------------------------------------------------------------------
# Row 24
# ---------------------------------------------------------------------------
$ws.Cells.Item(24, 3).Value = 35                      # C24
$ws.Cells.Item(24, 4).Value = 22                      # D24
$ws.Cells.Item(24, 5).Value = 59.090909090909         # E24
$ws.Cells.Item(24, 6).Value = 99                      # F24
$ws.Cells.Item(24, 7).Value = 95                      # G24
$ws.Cells.Item(24, 8).Value = 4.210526315789          # H24
$ws.Cells.Item(24, 9).Value = 924                     # I24
$ws.Cells.Item(24, 10).Value = 890                    # J24
$ws.Cells.Item(24, 11).Value = 3.820224719101         # K24
$ws.Cells.Item(24, 12).Value = 2.212389380530         # L24
$ws.Cells.Item(24, 13).Value = 41.500765696784        # M24

# ---------------------------------------------------------------------------
# Row 25
# ---------------------------------------------------------------------------
$ws.Cells.Item(25, 3).Value = 7                       # C25
$ws.Cells.Item(25, 4).Value = 3                       # D25
$ws.Cells.Item(25, 5).Value = 133.333333333333        # E25
$ws.Cells.Item(25, 6).Value = 32                      # F25
$ws.Cells.Item(25, 8).Value = 45.454545454545         # H25
$ws.Cells.Item(25, 9).Value = 264                     # I25
$ws.Cells.Item(25, 10).Value = 337                    # J25
$ws.Cells.Item(25, 11).Value = -21.661721068249       # K25
$ws.Cells.Item(25, 12).Value = -26.050420168067       # L25

# ---------------------------------------------------------------------------
# Row 26
# ---------------------------------------------------------------------------
$ws.Cells.Item(26, 3).Value = 4                       # C26
$ws.Cells.Item(26, 4).Value = 8                       # D26
$ws.Cells.Item(26, 5).Value = -50                     # E26
$ws.Cells.Item(26, 6).Value = 30                      # F26
$ws.Cells.Item(26, 7).Value = 35                      # G26
$ws.Cells.Item(26, 8).Value = -14.285714285714        # H26
$ws.Cells.Item(26, 9).Value = 306                     # I26
$ws.Cells.Item(26, 10).Value = 282                    # J26
$ws.Cells.Item(26, 11).Value = 8.510638297872         # K26
$ws.Cells.Item(26, 12).Value = 19.53125               # L26
$ws.Cells.Item(26, 13).Value = -4.672897196261        # M26

# ---------------------------------------------------------------------------
# Row 27
# ---------------------------------------------------------------------------
Set-TextZero $ws.Cells.Item(27, 6)     # F27 -> text "0"
$ws.Cells.Item(27, 7).Value = 1                       # G27
$ws.Cells.Item(27, 8).Value = -100                    # H27

# ---------------------------------------------------------------------------
# Row 28
# ---------------------------------------------------------------------------
$ws.Cells.Item(28, 6).Value = 1                       # F28
$ws.Cells.Item(28, 7).Value = 4                       # G28
$ws.Cells.Item(28, 8).Value = -75                     # H28
$ws.Cells.Item(28, 10).Value = 40                     # J28
$ws.Cells.Item(28, 11).Value = -17.5                  # K28
$ws.Cells.Item(28, 12).Value = 22.222222222222        # L28
